$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 3 (PORT 2) with the Inertial sensor entry
$ws.Range("B3").Value = "Inertial"
$ws.Range("C3").Value = "imu"
$ws.Range("D3").Value = "Inertial sensor for odometry"

# Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("F7").Select()
